# Fix things like missing parentheses and typos in English
# (and the corresponding translated strings).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the shared-string text content -----------------------------------
# Read the whole data grid (A1:F96) into memory so we can patch a handful
# of cells by their text content without disturbing anything else.
$rng = $ws.Range("A1:F96")
$arr = $rng.Value()

$nbsp = [char]0x00A0

# Row 47, column A (English): fix the "burried" -> "buried" typo.
$arr[47,1] = $arr[47,1].Replace("burried", "buried")

# Row 61: the photo-credit caption in each language was missing its
# closing parenthesis after "The Seattle Times". Add it back.
$arr[61,1] = $arr[61,1] + ")"   # English
$arr[61,2] = $arr[61,2] + ")"   # Chinese
# Column C (Russian) already had the closing parenthesis - leave as-is.
$arr[61,4] = $arr[61,4] + ")"   # Somali

# Column E (Spanish) also had non-breaking spaces around the slash;
# normalize them to regular spaces in addition to adding ")".
$arr[61,5] = $arr[61,5].Replace($nbsp, " ") + ")"   # Spanish

$arr[61,6] = $arr[61,6] + ")"   # Vietnamese

$rng.Value = $arr
